$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item('Citywide Totals')
$ws.Range("D2").Value = 91
$ws.Range("E2").Value = 74
$ws.Range("G2").Value = 85
$ws.Range("D3").Value = 138
$ws.Range("F3").Value = 141
$ws.Range("G3").Value = 142
$ws.Range("B9").Value = 382
$ws.Range("C9").Value = 485
$ws.Range("D9").Value = 418
$ws.Range("E9").Value = 482
$ws.Range("G9").Value = 438
$ws.Range("H9").Value = 445
$ws.Range("J9").Value = 424
$ws.Range("B10").Value = 1365
$ws.Range("C10").Value = 1614
$ws.Range("D10").Value = 1818
$ws.Range("E10").Value = 2218
$ws.Range("F10").Value = 2136
$ws.Range("G10").Value = 900
$ws.Range("H10").Value = 611
$ws.Range("I10").Value = 857
$ws.Range("J10").Value = 741
$ws.Range("B11").Value = 1883
$ws.Range("C11").Value = 2260
$ws.Range("D11").Value = 2478
$ws.Range("E11").Value = 2934
$ws.Range("F11").Value = 2935
$ws.Range("G11").Value = 1573
$ws.Range("H11").Value = 1340
$ws.Range("I11").Value = 1703
$ws.Range("J11").Value = 1554

$ws = $wb.Worksheets.Item('Chinatown')
$ws.Range("G3").Value = 1
$ws.Range("G9").Value = 14

$ws = $wb.Worksheets.Item('Garfield Park')
$ws.Range("D3").Value = 8
$ws.Range("D9").Value = 93

$ws = $wb.Worksheets.Item('Chatham')
$ws.Range("D8").Value = 16
$ws.Range("D10").Value = 59

$ws = $wb.Worksheets.Item('Loop')
$ws.Range("F3").Value = 14
$ws.Range("C7").Value = 41
$ws.Range("G7").Value = 56
$ws.Range("H7").Value = 64
$ws.Range("C8").Value = 334
$ws.Range("D8").Value = 518
$ws.Range("E8").Value = 670
$ws.Range("C9").Value = 393
$ws.Range("D9").Value = 590
$ws.Range("E9").Value = 756
$ws.Range("F9").Value = 635
$ws.Range("G9").Value = 245
$ws.Range("H9").Value = 207

$ws = $wb.Worksheets.Item('Armour Square')
$ws.Range("G6").Value = 11
$ws.Range("G7").Value = 15

$ws = $wb.Worksheets.Item('Old Town')
$ws.Range("J6").Value = 14
$ws.Range("J7").Value = 27

$ws = $wb.Worksheets.Item('Little Italy, UIC')
$ws.Range("J5").Value = 11
$ws.Range("J7").Value = 35

$ws = $wb.Worksheets.Item('North Lawndale')
$ws.Range("D7").Value = 49
$ws.Range("D8").Value = 76

$ws = $wb.Worksheets.Item('By Neighborhood')
$ws.Range("G5").Value = 15
$ws.Range("D8").Value = 62
$ws.Range("E8").Value = 119
$ws.Range("F8").Value = 158
$ws.Range("I8").Value = 83
$ws.Range("D19").Value = 59
$ws.Range("G21").Value = 14
$ws.Range("B27").Value = 24
$ws.Range("D27").Value = 27
$ws.Range("E27").Value = 25
$ws.Range("D28").Value = 102
$ws.Range("G28").Value = 85
$ws.Range("F30").Value = 15
$ws.Range("D32").Value = 93
$ws.Range("E41").Value = 27
$ws.Range("F41").Value = 33
$ws.Range("H41").Value = 18
$ws.Range("J50").Value = 35
$ws.Range("C53").Value = 393
$ws.Range("D53").Value = 590
$ws.Range("E53").Value = 756
$ws.Range("F53").Value = 635
$ws.Range("G53").Value = 245
$ws.Range("H53").Value = 207
$ws.Range("D54").Value = 19
$ws.Range("D61").Value = 27
$ws.Range("E61").Value = 62
$ws.Range("C62").Value = 29
$ws.Range("D65").Value = 76
$ws.Range("G66").Value = 3
$ws.Range("J70").Value = 27
$ws.Range("E76").Value = 98
$ws.Range("C78").Value = 36
$ws.Range("B89").Value = 26
$ws.Range("F89").Value = 21
$ws.Range("F98").Value = 13
$ws.Range("H98").Value = 8
$ws.Range("B99").Value = 1883
$ws.Range("C99").Value = 2260
$ws.Range("D99").Value = 2478
$ws.Range("E99").Value = 2934
$ws.Range("F99").Value = 2935
$ws.Range("G99").Value = 1573
$ws.Range("H99").Value = 1340
$ws.Range("I99").Value = 1703
$ws.Range("J99").Value = 1554

$ws = $wb.Worksheets.Item('Washington Park')
$ws.Range("F3").Value = 1
$ws.Range("B4").Value = 7
$ws.Range("B6").Value = 26
$ws.Range("F6").Value = 21

$ws = $wb.Worksheets.Item('Humboldt Park')
$ws.Range("E4").Value = 7
$ws.Range("F5").Value = 27
$ws.Range("H5").Value = 12
$ws.Range("E6").Value = 27
$ws.Range("F6").Value = 33
$ws.Range("H6").Value = 18

$ws = $wb.Worksheets.Item('Rush & Division')
$ws.Range("C5").Value = 32
$ws.Range("C6").Value = 36

$ws = $wb.Worksheets.Item('Englewood')
$ws.Range("D2").Value = 6
$ws.Range("G7").Value = 30
$ws.Range("D9").Value = 102
$ws.Range("G9").Value = 85

$ws = $wb.Worksheets.Item('Rogers Park')
$ws.Range("E9").Value = 79
$ws.Range("E10").Value = 98

$ws = $wb.Worksheets.Item('Edgewater')
$ws.Range("D3").Value = 1
$ws.Range("B6").Value = 20
$ws.Range("E6").Value = 19
$ws.Range("B7").Value = 24
$ws.Range("D7").Value = 27
$ws.Range("E7").Value = 25

$ws = $wb.Worksheets.Item('Near South Side')
$ws.Range("C6").Value = 26
$ws.Range("C7").Value = 29

$ws = $wb.Worksheets.Item('Lower West Side')
$ws.Range("D5").Value = 13
$ws.Range("D6").Value = 19

$ws = $wb.Worksheets.Item('Wrigleyville')
$ws.Range("H5").Value = 2
$ws.Range("F6").Value = 6
$ws.Range("F7").Value = 13
$ws.Range("H7").Value = 8

$ws = $wb.Worksheets.Item('Gage Park')
$ws.Range("F7").Value = 10
$ws.Range("F8").Value = 15

$ws = $wb.Worksheets.Item('Austin')
$ws.Range("E2").Value = 4
$ws.Range("D6").Value = 22
$ws.Range("F7").Value = 106
$ws.Range("I7").Value = 41
$ws.Range("D8").Value = 62
$ws.Range("E8").Value = 119
$ws.Range("F8").Value = 158
$ws.Range("I8").Value = 83

$ws = $wb.Worksheets.Item('North Park')
$ws.Range("G2").Value = 1
$ws.Range("G6").Value = 3
